# Update EPI Cartography 2; add license info
# Adds a new "module 3" block of vocabulary terms (confidentiality / ethics)
# to the vocabulary worksheet, appending rows 18-23 below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: module number for the new rows -----------------------------
for ($r = 18; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = 3
}

# --- Column B: terms --------------------------------------------------------
# Written in the order the terms should be registered so that the newly
# introduced vocabulary reads naturally top to bottom.
$ws.Cells.Item(18, 2).Value = "Confidentiality"
$ws.Cells.Item(19, 2).Value = "Privacy"
$ws.Cells.Item(21, 2).Value = "Belmont principles: justice"
$ws.Cells.Item(20, 2).Value = "Belmont principles: respect for persons"
$ws.Cells.Item(22, 2).Value = "Belmont principles: beneficence"
$ws.Cells.Item(23, 2).Value = "Geomask"

# --- Column C: definitions ---------------------------------------------------
$ws.Cells.Item(20, 3).Value = "Defined by two ethical convictions: a) individuals should be treated as autonomous agents; b) persons with diminished autonomy are entitled to protection"
$ws.Cells.Item(21, 3).Value = "Ethical principle that the burdens and benefits of research and public health practice should be justly distributed, including attention to need, effort, contribution, and merit"
$ws.Cells.Item(19, 3).Value = "The right of an individual to keep his or her information (health related or otherwise) private"
$ws.Cells.Item(18, 3).Value = "The duty of anyone entrusted with health information to keep that information private"
$ws.Cells.Item(22, 3).Value = "Two general rules have been formulated as complementary expressions of beneficent actions in this sense: (1) do not harm and (2) maximize possible benefits and minimize possible harms"
$ws.Cells.Item(23, 3).Value = "A class of methods for changing the geographic location of an individual in an unpredictable way to protect confidentiality, while trying to preserve the relationship between geocoded locations and disease occurrence (Sherman and Fetters 2007, Wiggins 2002)"

# Update the selected cell to reflect the position after the appended rows.
$ws.Application.Goto($ws.Range("A24"))
